$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that were tagged "line" / Line_1 now become "axis" / Axis_2
$line1Rows = 5,6,7,8,9,10,17,18,19,20

foreach ($r in $line1Rows) {
    $ws.Range("B$r").Value = "axis"
    $ws.Range("C$r").Value = "Axis_2"
    $ws.Range("D$r").ClearContents()
}

# Rows that were tagged "line" / Line_3 now become "axis" / Axis_3
$line3Rows = 14,15,21

foreach ($r in $line3Rows) {
    $ws.Range("B$r").Value = "axis"
    $ws.Range("C$r").Value = "Axis_3"
    $ws.Range("D$r").ClearContents()
}

# Update the active selection on the sheet
$null = $ws.Range("D8").Select()
